$wb = $excel.ActiveWorkbook

# ---- Sheet: VENTAS POR GRUPO ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M8").Value = 1653.76
$ws1.Range("I13").Value = 282.6
$ws1.Range("M13").Value = 5087.87
$ws1.Range("D21").Value = 979.58
$ws1.Range("D22").Value = "3 de 20"
$ws1.Range("I22").Value = "4 de 20"
$ws1.Range("M22").Value = "11 de 20"

# ---- Sheet: VENTA MENSUAL ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F8").Value = 2279.62
$ws2.Range("F13").Value = 5370.47
$ws2.Range("F21").Value = 5388.82
$ws2.Range("F22").Value = 57560.41

# ---- Sheet: CUMPLIMIENTO MENSUAL ----
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Excel COM stores column width with a +5/MaximumDigitWidth (=0.8333 for
# Calibri 11, MDW=6px) padding baked into the serialized <col width=.../>
# value, so back the entered ColumnWidth off by that amount to land the
# saved OOXML width attribute exactly on 24 (matching the target diff).
$ws3.Columns.Item(5).ColumnWidth = 23.166666666666668

$ws3.Range("D3").Value = 3894.62
$ws3.Range("E3").Value = 273.4515657367901
$ws3.Range("F3").Value = 0.9343937450631435

$ws3.Range("D8").Value = 1470.6
$ws3.Range("E8").Value = -845.5999999999999
$ws3.Range("F8").Value = 2.35296

$ws3.Range("D16").Value = 46464.18
$ws3.Range("E16").Value = -2197.940000000002
$ws3.Range("F16").Value = 1.049652737616748

$ws3.Range("D19").Value = 57560.41
$ws3.Range("E19").Value = 7817.587622917684
$ws3.Range("F19").Value = 0.8804247926342531
